$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.984.30'
$ws.Range("E2").Value = '  -5.05%  '

$ws.Range("D3").Value = '2.549.63'
$ws.Range("E3").Value = '  -5.64%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = "'299.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.35%  '

$ws.Range("D6").Value = "'94.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.34%  '

$ws.Range("E7").Value = '  -4.23%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").Value = "'0.550"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.79%  '

$ws.Range("D10").Value = "'36.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.66%  '

$ws.Range("D11").Value = "'0.0810"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.46%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = "'7.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.43%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = "'0.115"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.49%  '

$ws.Range("D14").Value = '2.939.35'
$ws.Range("E14").Value = '  -5.72%  '

$ws.Range("D15").Value = '2.540.16'
$ws.Range("E15").Value = '  -6.10%  '

$ws.Range("D16").Value = "'0.881"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.69%  '

$ws.Range("D17").Value = "'14.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.52%  '

$ws.Range("D18").Value = '42.990.05'
$ws.Range("E18").Value = '  -5.27%  '

$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").Value = "'12.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.55%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0982'
$ws.Range("E20").Value = '  -3.41%  '

$ws.Range("D21").Value = "'6.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.77%  '

$ws.Range("D22").Value = "'72.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.68%  '

$ws.Range("D23").Value = "'254.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -9.44%  '

$ws.Range("D24").Value = "'2.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.27%  '

$ws.Range("D25").Value = "'2.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.76%  '

$ws.Range("D26").Value = "'29.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.76%  '

$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("D28").Value = "'10.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.59%  '

$ws.Range("D29").Value = "'36.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.65%  '

$ws.Range("D30").Value = "'2.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.08%  '

$ws.Range("D31").Value = "'6.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.28%  '

$ws.Range("D32").Value = "'152.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.81%  '

$ws.Range("E33").Value = '  -1.94%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = "'3.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -11.67%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = "'2.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -10.47%  '

$ws.Range("D36").Value = "'0.0793"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.34%  '

$ws.Range("E37").Value = '  -6.23%  '

$ws.Range("D38").Value = "'17.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.73%  '

$ws.Range("E39").Value = '  -4.25%  '

$ws.Range("D40").Value = "'22.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -11.26%  '

$ws.Range("D41").Value = "'3.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.07%  '

$ws.Range("D42").Value = "'0.0310"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.12%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = "'3.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.43%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.109.40'
$ws.Range("E44").Value = '  -1.39%  '

$ws.Range("D45").Value = "'1.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +20.67%  '

$ws.Range("D46").Value = "'0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.16%  '

$ws.Range("E47").Value = '  -4.33%  '

$ws.Range("D48").Value = "'84.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -9.86%  '

$ws.Range("D49").Value = '2.792.44'
$ws.Range("E49").Value = '  -5.67%  '

$ws.Range("D50").Value = "'104.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.99%  '

$ws.Range("D51").Value = "'1.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.57%  '
